$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '29.713.89'
$ws.Range("E2").Value = '  -2.61%  '

$ws.Range("D3").Value = '2.094.37'
$ws.Range("E3").Value = '  -1.95%  '

$ws.Range("E4").Value = '  +0.45%  '

$cell = $ws.Range("D5")
$cell.NumberFormat = "@"
$cell.Value = '343.27'
$cell.Style = "Normal"
$ws.Range("E5").Value = '  -2.43%  '

$ws.Range("E6").Value = '  +0.35%  '

$cell = $ws.Range("D7")
$cell.NumberFormat = "@"
$cell.Value = '0.5166'
$cell.Style = "Normal"
$ws.Range("E7").Value = '  -1.62%  '

$cell = $ws.Range("D8")
$cell.NumberFormat = "@"
$cell.Value = '0.4380'
$cell.Style = "Normal"
$ws.Range("E8").Value = '  -3.72%  '

$cell = $ws.Range("D9")
$cell.NumberFormat = "@"
$cell.Value = '52.87'
$cell.Style = "Normal"
$ws.Range("E9").Value = '  -1.36%  '

$cell = $ws.Range("D10")
$cell.NumberFormat = "@"
$cell.Value = '0.09276'
$cell.Style = "Normal"
$ws.Range("E10").Value = '  +1.35%  '

$ws.Range("E11").Value = '  -2.27%  '

$cell = $ws.Range("D12")
$cell.NumberFormat = "@"
$cell.Value = '24.88'
$cell.Style = "Normal"
$ws.Range("E12").Value = '  -2.00%  '

$ws.Range("D13").Value = '2.104.03'
$ws.Range("E13").Value = '  -1.72%  '

$cell = $ws.Range("D14")
$cell.NumberFormat = "@"
$cell.Value = '8.270'
$cell.Style = "Normal"
$ws.Range("E14").Value = '  +1.42%  '

$cell = $ws.Range("D15")
$cell.NumberFormat = "@"
$cell.Value = '6.750'
$cell.Style = "Normal"
$ws.Range("E15").Value = '  -1.95%  '

$cell = $ws.Range("D16")
$cell.NumberFormat = "@"
$cell.Value = '99.47'
$cell.Style = "Normal"
$ws.Range("E16").Value = '  -2.33%  '

$cell = $ws.Range("D17")
$cell.NumberFormat = "@"
$cell.Value = '0.00001154'
$cell.Style = "Normal"
$ws.Range("E17").Value = '  -1.30%  '

$cell = $ws.Range("D18")
$cell.NumberFormat = "@"
$cell.Value = '1.010'
$cell.Style = "Normal"
$ws.Range("E18").Value = '  +0.27%  '

$cell = $ws.Range("D19")
$cell.NumberFormat = "@"
$cell.Value = '20.76'
$cell.Style = "Normal"
$ws.Range("E19").Value = '  +1.75%  '

$cell = $ws.Range("D20")
$cell.NumberFormat = "@"
$cell.Value = '0.06640'
$cell.Style = "Normal"
$ws.Range("E20").Value = '  -1.18%  '

$ws.Range("E21").Value = '  +0.30%  '

$cell = $ws.Range("D22")
$cell.NumberFormat = "@"
$cell.Value = '6.190'
$cell.Style = "Normal"
$ws.Range("E22").Value = '  -2.84%  '

$ws.Range("D23").Value = '29.745.66'
$ws.Range("E23").Value = '  -2.77%  '

$cell = $ws.Range("D24")
$cell.NumberFormat = "@"
$cell.Value = '12.49'
$cell.Style = "Normal"
$ws.Range("E24").Value = '  -2.91%  '

$cell = $ws.Range("D25")
$cell.NumberFormat = "@"
$cell.Value = '2.321'
$cell.Style = "Normal"
$ws.Range("E25").Value = '  -2.65%  '

$ws.Range("D26").Value = '2.352.93'
$ws.Range("E26").Value = '  -1.61%  '

$cell = $ws.Range("D27")
$cell.NumberFormat = "@"
$cell.Value = '21.93'
$cell.Style = "Normal"
$ws.Range("E27").Value = '  -2.58%  '

$cell = $ws.Range("D28")
$cell.NumberFormat = "@"
$cell.Value = '2.510'
$cell.Style = "Normal"
$ws.Range("E28").Value = '  -4.12%  '

$cell = $ws.Range("D29")
$cell.NumberFormat = "@"
$cell.Value = '161.23'
$cell.Style = "Normal"
$ws.Range("E29").Value = '  -2.07%  '

$cell = $ws.Range("D30")
$cell.NumberFormat = "@"
$cell.Value = '132.96'
$cell.Style = "Normal"
$ws.Range("E30").Value = '  -2.12%  '

$cell = $ws.Range("D31")
$cell.NumberFormat = "@"
$cell.Value = '1.135'
$cell.Style = "Normal"
$ws.Range("E31").Value = '  -7.38%  '

$ws.Range("E32").Value = '  -2.96%  '

$cell = $ws.Range("D33")
$cell.NumberFormat = "@"
$cell.Value = '1.651'
$cell.Style = "Normal"
$ws.Range("E33").Value = '  -4.26%  '

$cell = $ws.Range("D34")
$cell.NumberFormat = "@"
$cell.Value = '6.155'
$cell.Style = "Normal"
$ws.Range("E34").Value = '  -3.61%  '

$cell = $ws.Range("D35")
$cell.NumberFormat = "@"
$cell.Value = '3.940'
$cell.Style = "Normal"
$ws.Range("E35").Value = '  -2.21%  '

$cell = $ws.Range("D36")
$cell.NumberFormat = "@"
$cell.Value = '6.260'
$cell.Style = "Normal"
$ws.Range("E36").Value = '  +2.16%  '

$cell = $ws.Range("D37")
$cell.NumberFormat = "@"
$cell.Value = '10.22'
$cell.Style = "Normal"
$ws.Range("E37").Value = '  -2.19%  '

$cell = $ws.Range("D38")
$cell.NumberFormat = "@"
$cell.Value = '0.02574'
$cell.Style = "Normal"
$ws.Range("E38").Value = '  -2.64%  '

$cell = $ws.Range("D39")
$cell.NumberFormat = "@"
$cell.Value = '0.06709'
$cell.Style = "Normal"
$ws.Range("E39").Value = '  -3.76%  '

$cell = $ws.Range("D40")
$cell.NumberFormat = "@"
$cell.Value = '12.46'
$cell.Style = "Normal"
$ws.Range("E40").Value = '  -1.93%  '

$ws.Range("E41").Value = '  -1.16%  '

$ws.Range("E42").Value = '  -5.09%  '

$cell = $ws.Range("D43")
$cell.NumberFormat = "@"
$cell.Value = '1.315'
$cell.Style = "Normal"
$ws.Range("E43").Value = '  +3.05%  '

$cell = $ws.Range("D44")
$cell.NumberFormat = "@"
$cell.Value = '0.6771'
$cell.Style = "Normal"
$ws.Range("E44").Value = '  +4.01%  '

$cell = $ws.Range("D45")
$cell.NumberFormat = "@"
$cell.Value = '14.29'
$cell.Style = "Normal"
$ws.Range("E45").Value = '  -3.69%  '

$cell = $ws.Range("D46")
$cell.NumberFormat = "@"
$cell.Value = '2.319'
$cell.Style = "Normal"
$ws.Range("E46").Value = '  -1.30%  '

$cell = $ws.Range("D47")
$cell.NumberFormat = "@"
$cell.Value = '0.00000000362'
$cell.Style = "Normal"
$ws.Range("E47").Value = '  -3.36%  '

$cell = $ws.Range("D48")
$cell.NumberFormat = "@"
$cell.Value = '3.619'
$cell.Style = "Normal"
$ws.Range("E48").Value = '  -3.32%  '

$ws.Range("E49").Value = '  -2.51%  '

$cell = $ws.Range("D50")
$cell.NumberFormat = "@"
$cell.Value = '81.95'
$cell.Style = "Normal"
$ws.Range("E50").Value = '  -2.25%  '

$cell = $ws.Range("D51")
$cell.NumberFormat = "@"
$cell.Value = '1.161'
$cell.Style = "Normal"
$ws.Range("E51").Value = '  -1.97%  '
